$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Widen column E to fit the new (longer) answer text
$ws.Columns.Item(5).ColumnWidth = 28

# --- New question block (nr 6): "Czy gusu przyjmie projek?" ---
# Column D (question text) is identical for all 6 rows.
$ws.Cells.Item(24, 4).Value = "Czy gusu przyjmie projek?"
$ws.Cells.Item(25, 4).Value = "Czy gusu przyjmie projek?"
$ws.Cells.Item(26, 4).Value = "Czy gusu przyjmie projek?"
$ws.Cells.Item(27, 4).Value = "Czy gusu przyjmie projek?"
$ws.Cells.Item(28, 4).Value = "Czy gusu przyjmie projek?"
$ws.Cells.Item(29, 4).Value = "Czy gusu przyjmie projek?"

# Column E (answer text) - entered in the order: wrong, reused "NIE WIEM", correct,
# correct alt, then wrong extras added afterwards.
$ws.Cells.Item(24, 5).Value = "uep lenie"
$ws.Cells.Item(25, 5).Value = "NIE WIEM"
$ws.Cells.Item(26, 5).Value = "tak "
$ws.Cells.Item(29, 5).Value = "oczywiście estymator krul"

# Column H (image filename) - typed once then copied down.
$ws.Cells.Item(24, 8).Value = "2.jpg"
$ws.Cells.Item(25, 8).Value = "2.jpg"
$ws.Cells.Item(26, 8).Value = "2.jpg"
$ws.Cells.Item(27, 8).Value = "2.jpg"
$ws.Cells.Item(28, 8).Value = "2.jpg"
$ws.Cells.Item(29, 8).Value = "2.jpg"

# Remaining two wrong-answer rows added afterwards.
$ws.Cells.Item(27, 5).Value = "chyba"
$ws.Cells.Item(28, 5).Value = "może"

# Column A (question number)
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(26, 1).Value = 6
$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(28, 1).Value = 6
$ws.Cells.Item(29, 1).Value = 6

# Column B (category)
$ws.Cells.Item(24, 2).Value = "ge"
$ws.Cells.Item(25, 2).Value = "ge"
$ws.Cells.Item(26, 2).Value = "ge"
$ws.Cells.Item(27, 2).Value = "ge"
$ws.Cells.Item(28, 2).Value = "ge"
$ws.Cells.Item(29, 2).Value = "ge"

# Column C (difficulty)
$ws.Cells.Item(24, 3).Value = 2
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(28, 3).Value = 2
$ws.Cells.Item(29, 3).Value = 2

# Column F (odp - is this row the displayed/"given" answer flag)
$ws.Cells.Item(24, 6).Value = "F"
$ws.Cells.Item(25, 6).Value = "F"
$ws.Cells.Item(26, 6).Value = "F"
$ws.Cells.Item(27, 6).Value = "F"
$ws.Cells.Item(28, 6).Value = "F"
$ws.Cells.Item(29, 6).Value = "T"

# Column G (praw - correct answer flag)
$ws.Cells.Item(24, 7).Value = "T"
$ws.Cells.Item(25, 7).Value = "T"
$ws.Cells.Item(26, 7).Value = "T"
$ws.Cells.Item(27, 7).Value = "T"
$ws.Cells.Item(28, 7).Value = "T"
$ws.Cells.Item(29, 7).Value = "T"

# Update the view: this is the panel shown once all quiz questions are finished.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E31").Select()
